# Update computed power-flow results (pl_mw) for the 380 kV case
# Columns updated per row: B,C,D,F,G,J,K,L,M,N (column indexes below)
# Rows 2-25 correspond to data rows (A2:A25 holds the 0-based index column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2,3,4,6,7,10,11,12,13,14)

$rowData = @{
  2 = @(2.035694872141079,0.0320291313009875,0.007403928998831688,4.527041862629872,0.002662128419459988,0.2025864728366589,1.476498679640741,0.2969469157344946,0.4424505019053129,4.805848590608093)
  3 = @(2.005020316188791,0.02787730805349042,0.007423560904857851,4.517611846578887,0.002666457437551046,0.2031645383232572,1.443347790476111,0.296313459523283,0.4378898275139242,4.813582738122179)
  4 = @(1.987258056098369,0.02533135446243762,0.007439657779967135,4.51332765514924,0.002669257677515449,0.2035509152456498,1.423838058129604,0.2960433045902207,0.4353074513802397,4.819218206949401)
  5 = @(1.980289621505818,0.02429462227439672,0.007447241191099963,4.511960660552077,0.002670434670860079,0.2037162895117071,1.416100252027917,0.2959631465779182,0.4343099765852578,4.821737668984241)
  6 = @(1.979148823110762,0.02412251809255395,0.007448562460497854,4.511756559205764,0.002670632279426954,0.2037442287579392,1.414828236188953,0.2959516460652836,0.4341476625023866,4.822169490277361)
  7 = @(1.987162984476385,0.02531736971347698,0.007439755897584277,4.513307685210918,0.002669273405355362,0.2035531134416555,1.423732842589203,0.296042102273546,0.4352937768454481,4.819251282546816)
  8 = @(2.024895920392623,0.03059686342962209,0.007409863253360349,4.523477945511118,0.002663591618729224,0.2027792735776259,1.464892891666977,0.2967038812149738,0.4408327976415869,4.808331293837526)
  9 = @(2.107391729197502,0.0409792241174074,0.007383018766535798,4.555367420138737,0.002653572789122065,0.2015105564555135,1.552315972726234,0.2989418464846096,0.4534216734663588,4.793953918302208)
  10 = @(2.173188169605112,0.04863020533551321,0.007382279694958527,4.586082947470288,0.002646889367872603,0.200729160615662,1.620649399530976,0.3011567071208106,0.4637224405779605,4.787684591540895)
  11 = @(2.204248543371023,0.05211700840814615,0.007385987365736923,4.601639893361096,0.002643994449450487,0.2004062185596851,1.652630851347624,0.3022877172720584,0.4686368137690522,4.785765941186327)
  12 = @(2.216172621852877,0.05343835938522545,0.007387966470796314,4.607758717016651,0.002642919010080929,0.2002885894888529,1.664870389594967,0.3027337041704499,0.4705305782561595,4.785173674744527)
  13 = @(2.21359734782402,0.05315373846022453,0.007387514758269376,4.606430790899367,0.002643149701690863,0.200313715889525,1.662228655372161,0.3026368666677826,0.4701212644361377,4.785295256274907)
  14 = @(2.205226295063937,0.05222569685528811,0.007386138697679456,4.602138728698492,0.002643905556095315,0.200396447790645,1.653635223144732,0.3023240545047443,0.4687919579330568,4.785714523342477)
  15 = @(2.200119899682363,0.05165737335691745,0.007385370528465884,4.599539373829515,0.002644371244830079,0.2004477301927636,1.648388278193352,0.3021347511805885,0.4679819888183374,4.785988826753169)
  16 = @(2.171181007577559,0.04840246742503496,0.007382118211733513,4.585098138075239,0.002647081474494923,0.2007509187929806,1.618577369029282,0.3010852725104556,0.4634058675957675,4.787828766731195)
  17 = @(2.153717019997316,0.04640735530870188,0.007381154136316681,4.576644629246218,0.002648781278671033,0.200945233962301,1.600518871919775,0.3004730325449074,0.4606570470430214,4.789196604186074)
  18 = @(2.143778497123378,0.04526041188825047,0.007380980627888611,4.571931520121808,0.002649772653393414,0.2010600605687429,1.590216460563084,0.3001325121176848,0.4590975094301228,4.790071189512147)
  19 = @(2.14043174738066,0.04487217634930118,0.007380987516866888,4.570361357226702,0.002650110670918609,0.2010994651594942,1.586742730497889,0.3000192158491828,0.4585731732342921,4.790382393653388)
  20 = @(2.155565090977461,0.04661967664060285,0.007381217378956251,4.577529086301212,0.002648598915455701,0.2009242320133033,1.602432499141486,0.3005370039172348,0.4609474377896916,4.789041903726911)
  21 = @(2.207680675929169,0.05249825821883292,0.007386527322209346,4.603393231003224,0.002643682979447789,0.2003720210102671,1.656155825501685,0.3024154551413432,0.4691815178214611,4.785587729347412)
  22 = @(2.242686500021193,0.05634596004475156,0.007393346897282882,4.621624393042509,0.002640591341769083,0.2000382869042632,1.692018194362163,0.3037462633042196,0.4747541100189352,4.784112980625707)
  23 = @(2.223916802478527,0.05429182567648638,0.007389402781948107,4.611772639702963,0.002642230351092405,0.2002139258514752,1.672809055446862,0.3030265674969073,0.4717624425054865,4.784828432923064)
  24 = @(2.154729261106468,0.04652368594484813,0.007381187600415728,4.577128765566755,0.002648681317770103,0.2009337172986747,1.601567100140272,0.3005080467455699,0.4608160873506293,4.789111569038539)
  25 = @(2.08416388757351,0.03816684817384441,0.007386919869417419,4.545461230585119,0.002656163664111299,0.2018272412963178,1.527946103996157,0.2982360272937257,0.4498312867634624,4.797089659075596)
}

foreach ($rowNum in $rowData.Keys) {
  $values = $rowData[$rowNum]
  for ($i = 0; $i -lt $colIndexes.Length; $i++) {
    $ws.Cells.Item([int]$rowNum, [int]$colIndexes[$i]).Value = [double]$values[$i]
  }
}

Write-Output "Updated $($rowData.Keys.Count) rows x $($colIndexes.Length) columns"
